$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-05-17 Saturday" "2025-05-18 Sunday"

Replace-Text "96×21=2016" "91×59=5369"
Replace-Text "73×88=6424" "48×11=528"
Replace-Text "61×42=2562" "29×43=1247"
Replace-Text "56×79=4424" "27×83=2241"
Replace-Text "33×12=396" "50×65=3250"

Replace-Text "68×76=5168" "81×63=5103"
Replace-Text "62×45=2790" "47×61=2867"
Replace-Text "51×39=1989" "61×56=3416"
Replace-Text "39×80=3120" "94×44=4136"
Replace-Text "53×13=689" "81×80=6480"

Replace-Text "72×33=2376" "31×54=1674"
Replace-Text "83×59=4897" "87×12=1044"
Replace-Text "49×88=4312" "72×83=5976"
Replace-Text "66×31=2046" "21×79=1659"
Replace-Text "16×71=1136" "51×64=3264"

Replace-Text "21×74=1554" "67×11=737"
Replace-Text "88×64=5632" "47×26=1222"
Replace-Text "27×89=2403" "27×36=972"
Replace-Text "61×23=1403" "55×76=4180"
Replace-Text "84×34=2856" "75×12=900"

Replace-Text "88×13=1144" "65×20=1300"
Replace-Text "67×70=4690" "94×51=4794"
Replace-Text "88×84=7392" "31×39=1209"
Replace-Text "38×27=1026" "64×74=4736"
Replace-Text "33×45=1485" "97×43=4171"
